$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.462.03"
$ws.Range("E2").Value = "  +1.63%  "

# Row 3
$ws.Range("D3").Value = "2.163.48"
$ws.Range("E3").Value = "  +3.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.07"
$ws.Range("E5").Value = "  +0.22%  "

# Row 6
$ws.Range("E6").Value = "  +1.15%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.23"
$ws.Range("E7").Value = "  +4.90%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  +3.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0864"
$ws.Range("E10").Value = "  +2.59%  "

# Row 11
$ws.Range("E11").Value = "  -0.68%  "

# Row 12
$ws.Range("E12").Value = "  +7.05%  "

# Row 13
$ws.Range("D13").Value = "2.483.14"
$ws.Range("E13").Value = "  +3.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.25"
$ws.Range("E14").Value = "  +1.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.819"
$ws.Range("E15").Value = "  +2.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.58"
$ws.Range("E16").Value = "  +1.96%  "

# Row 17
$ws.Range("D17").Value = "2.162.31"
$ws.Range("E17").Value = "  +3.28%  "

# Row 18
$ws.Range("D18").Value = "39.424.75"
$ws.Range("E18").Value = "  +1.72%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.37"
$ws.Range("E19").Value = "  +1.06%  "

# Row 20
$ws.Range("E20").Value = "  +1.71%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  +1.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.55"
$ws.Range("E22").Value = "  +0.49%  "

# Row 23
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  -1.41%  "

# Row 25
$ws.Range("E25").Value = "  +1.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.67"
$ws.Range("E26").Value = "  +1.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.04"
$ws.Range("E27").Value = "  +0.58%  "

# Row 28
$ws.Range("E28").Value = "  -1.43%  "

# Row 29
$ws.Range("E29").Value = "  -3.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.67"
$ws.Range("E30").Value = "  +2.56%  "

# Row 31
$ws.Range("E31").Value = "  +8.23%  "

# Row 32
$ws.Range("E32").Value = "  +1.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.68"
$ws.Range("E33").Value = "  +3.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  +2.39%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.06"
$ws.Range("E35").Value = "  +8.75%  "

# Row 36
$ws.Range("E36").Value = "  +2.08%  "

# Row 37
$ws.Range("E37").Value = "  +2.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("E38").Value = "  -0.61%  "

# Row 39
$ws.Range("E39").Value = "  -0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.23"
$ws.Range("E40").Value = "  +0.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.76"
$ws.Range("E41").Value = "  +2.88%  "

# Row 42
$ws.Range("E42").Value = "  +2.12%  "

# Row 43
$ws.Range("D43").Value = "1.532.20"
$ws.Range("E43").Value = "  -0.70%  "

# Row 44
$ws.Range("E44").Value = "  +5.20%  "

# Row 45
$ws.Range("E45").Value = "  +1.07%  "

# Row 46
$ws.Range("E46").Value = "  +6.98%  "

# Row 48
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.24"
$ws.Range("E48").Value = "  +3.55%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.78"
$ws.Range("E49").Value = "  +1.05%  "

# Row 50
$ws.Range("D50").Value = "2.366.56"
$ws.Range("E50").Value = "  +3.34%  "

# Row 51
$ws.Range("E51").Value = "  +0.27%  "
